# Add "ejecutor" rows (EJecutor de tareas en automatico)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "29-01-2024"
$ws.Range("B12").Value = "17:12:09"
$ws.Range("C12").Value = "registro_test2"
$ws.Range("D12").Value = 2.53

$ws.Range("A13").Value = "29-01-2024"
$ws.Range("B13").Value = "17:12:15"
$ws.Range("C13").Value = "test_form2"
$ws.Range("D13").Value = 5.3
